$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A, E, H, I sometimes hold digit-only strings ("877", "09087011142",
# "300", "10", ...) that Excel's COM layer would otherwise coerce to numbers
# (and, for the mobile numbers, silently drop the leading zero). Force those
# ranges to Text format before writing so the values land as real strings,
# matching the rest of the sheet's shared-string-backed columns. Each range
# is handled separately (rather than as one multi-area range) since the
# number-format assignment only reliably reaches the first area otherwise.
$colA = $ws.Range("A7:A9")
$colE = $ws.Range("E7:E9")
$colH = $ws.Range("H7:H9")
$colI = $ws.Range("I7:I9")
$colA.NumberFormat = "@"
$colE.NumberFormat = "@"
$colH.NumberFormat = "@"
$colI.NumberFormat = "@"

# --- Row 7: invoice 11-004 / dhanush / fireworks ---
$ws.Range("A7").Value = "11-004"
$ws.Range("B7").Value = "23.11.2025"
$ws.Range("C7").Value = "dhanush"
$ws.Range("D7").Value = "Boys hostel Banari amman institute of technology"
$ws.Range("E7").Value = "09087011142"
$ws.Range("F7").Value = "dhanushh@gamil"
$ws.Range("G7").Value = "fireworks"
$ws.Range("H7").Value = "300"
$ws.Range("I7").Value = "10"
$ws.Range("J7").Value = 3000
$ws.Range("K7").Value = 300
$ws.Range("L7").Value = 200
$ws.Range("M7").Value = 3100
$ws.Range("N7").Value = 100
$ws.Range("O7").Value = 3000

# --- Row 8: invoice 877 / dhanush / fireworks ---
$ws.Range("A8").Value = "877"
$ws.Range("B8").Value = "23.11.2025"
$ws.Range("C8").Value = "dhanush"
$ws.Range("D8").Value = "Boys hostel Banari amman institute of technology"
$ws.Range("E8").Value = "09087011142"
$ws.Range("F8").Value = "dhanushh@gamil"
$ws.Range("G8").Value = "fireworks"
$ws.Range("H8").Value = "300"
$ws.Range("I8").Value = "10"
$ws.Range("J8").Value = 3000
$ws.Range("K8").Value = 300
$ws.Range("L8").Value = 200
$ws.Range("M8").Value = 3100
$ws.Range("N8").Value = 100
$ws.Range("O8").Value = 3000

# --- Row 9: invoice 432 / dhanush / fireworks ---
$ws.Range("A9").Value = "432"
$ws.Range("B9").Value = "23.11.2025"
$ws.Range("C9").Value = "dhanush"
$ws.Range("D9").Value = "Boys hostel Banari amman institute of technology"
$ws.Range("E9").Value = "09087011142"
$ws.Range("F9").Value = "dhanushh@gamil"
$ws.Range("G9").Value = "fireworks"
$ws.Range("H9").Value = "300"
$ws.Range("I9").Value = "10"
$ws.Range("J9").Value = 3000
$ws.Range("K9").Value = 300
$ws.Range("L9").Value = 200
$ws.Range("M9").Value = 3100
$ws.Range("N9").Value = 100
$ws.Range("O9").Value = 3000

# Drop the temporary text format now that the strings are committed so no
# extra cell formatting lingers on the sheet.
$colA.ClearFormats()
$colE.ClearFormats()
$colH.ClearFormats()
$colI.ClearFormats()

# --- Row 10 no longer exists in the updated invoice log: delete it so the
#     remaining data shifts up and the sheet ends at row 9. ---
$ws.Rows.Item(10).Delete()

# --- View state: zoom + active selection matching the edited workbook ---
$ws.Range("P5").Select()
$excel.ActiveWindow.Zoom = 100
